$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of data (row 6): research paper title, year, lipase source.
# The pre-existing D6 value ("Glycerol/Free fatty acids/...") is left untouched.
$ws.Range("A6").Value = "Pretreatment of Coconut Mill Effluent Using Celite-Immobilized Hydrolytic Enzyme Preparation from Staphylococcus pasteuri and Its Impact on Anaerobic Digestion"
$ws.Range("B6").Value = 2015
$ws.Range("C6").Value = "Lipase from Staphylococcus pasteuri"

# Row-height adjustments: row 3 shrinks (45 -> 30), row 6 grows (30 -> 60)
# now that it holds a full four-column record.
$ws.Rows.Item(3).RowHeight = 30
$ws.Rows.Item(6).RowHeight = 60

# Move the active selection to the newly completed row.
[void]$ws.Range("D6").Select()
